$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder country labels that changed rank (swap/rotate via captured text) ---

# Venezuela's total cases overtook Nepal's: row105 becomes Venezuela, row106 becomes Nepal
$a105 = $ws.Range("A105").Text
$a106 = $ws.Range("A106").Text
$ws.Range("A105").Value = $a106
$ws.Range("A106").Value = $a105

# Zimbabue's total cases overtook Gibraltar's and Guadalupe's:
# row162 becomes Zimbabue, row163 becomes Gibraltar, row164 becomes Guadalupe
$a162 = $ws.Range("A162").Text
$a163 = $ws.Range("A163").Text
$a164 = $ws.Range("A164").Text
$ws.Range("A162").Value = $a164
$ws.Range("A163").Value = $a162
$ws.Range("A164").Value = $a163

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados,
#     Casos criticos, Muertes hoy, Muertes) ---

# Estados Unidos
$ws.Range("B4").Value = 1815762
$ws.Range("C4").Value = 22232
$ws.Range("D4").Value = 530474
$ws.Range("E4").Value = 1179752
$ws.Range("G4").Value = 994
$ws.Range("H4").Value = 105536

# Brasil
$ws.Range("B5").Value = 498440
$ws.Range("C5").Value = 30102
$ws.Range("D5").Value = 205371
$ws.Range("E5").Value = 264235
$ws.Range("G5").Value = 890
$ws.Range("H5").Value = 28834

# Colombia
$ws.Range("B34").Value = 28236
$ws.Range("C34").Value = 1548
$ws.Range("D34").Value = 7121
$ws.Range("E34").Value = 20225
$ws.Range("G34").Value = 37
$ws.Range("H34").Value = 890

# Japon
$ws.Range("B45").Value = 16804
$ws.Range("C45").Value = 85
$ws.Range("D45").Value = 14406
$ws.Range("E45").Value = 1512
$ws.Range("G45").Value = 12
$ws.Range("H45").Value = 886

# Argentina
$ws.Range("B47").Value = 16214
$ws.Range("C47").Value = 795
$ws.Range("D47").Value = 4788
$ws.Range("E47").Value = 10898
$ws.Range("G47").Value = 8
$ws.Range("H47").Value = 528

# Nigeria
$ws.Range("B56").Value = 9855
$ws.Range("C56").Value = 553
$ws.Range("D56").Value = 2856
$ws.Range("E56").Value = 6726
$ws.Range("G56").Value = 12
$ws.Range("H56").Value = 273

# Chequia
$ws.Range("B58").Value = 9230
$ws.Range("C58").Value = 34
$ws.Range("D58").Value = 6546
$ws.Range("E58").Value = 2365

# Noruega
$ws.Range("B61").Value = 8437
$ws.Range("C61").Value = 15
$ws.Range("D61").Value = 7727
$ws.Range("E61").Value = 474

# Venezuela (now at row 105)
$ws.Range("B105").Value = 1459
$ws.Range("C105").Value = 89
$ws.Range("D105").Value = 302
$ws.Range("E105").Value = 1143
$ws.Range("H105").Value = 14

# Nepal (now at row 106)
$ws.Range("B106").Value = 1401
$ws.Range("C106").Value = 189
$ws.Range("D106").Value = 219
$ws.Range("E106").Value = 1176
$ws.Range("H106").Value = 6

# Uruguay
$ws.Range("B125").Value = 821
$ws.Range("C125").Value = 5
$ws.Range("D125").Value = 682
$ws.Range("E125").Value = 117

# Togo
$ws.Range("B144").Value = 433
$ws.Range("C144").Value = 5
$ws.Range("D144").Value = 206
$ws.Range("E144").Value = 214

# Zimbabue (now at row 162)
$ws.Range("B162").Value = 174
$ws.Range("C162").Value = 25
$ws.Range("D162").Value = 29
$ws.Range("E162").Value = 141
$ws.Range("H162").Value = 4

# Gibraltar (now at row 163)
$ws.Range("B163").Value = 169
$ws.Range("C163").Value = 8
$ws.Range("D163").Value = 149
$ws.Range("E163").Value = 20
$ws.Range("H163").Value = 0

# Guadalupe (now at row 164)
$ws.Range("B164").Value = 162
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 138
$ws.Range("E164").Value = 10
$ws.Range("H164").Value = 14

# --- Update "last updated" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 01:05"
